# add role for upload-image
# Inserts three new permission rows ("图片管理" / "追加图片" / "删除图片",
# grouped under "商品管理 - 图片管理") right before the existing "14000"
# ("商品内容管理") row, pushing everything from old row 35 onward down by 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "14000" entry currently lives on row 35 - insert three blank rows
# above it so the new rows land at 35, 36 and 37.
$ws.Rows.Item(35).Insert()
$ws.Rows.Item(36).Insert()
$ws.Rows.Item(37).Insert()

# Match the existing data-row formatting (font + thin box border, same as
# style index 3 used by every other data row) for the new A:C cells.
$fmt = $ws.Range("A35:C37")
$fmt.Font.Name = "微软雅黑"
$fmt.Font.Size = 12
$fmt.Borders.LineStyle = 1
$fmt.Borders.Weight = 2

# Fill in the IDs first …
$ws.Range("A35").Value = 13500
$ws.Range("A36").Value = 13502
$ws.Range("A37").Value = 13503

# … then the permission names …
$ws.Range("B35").Value = "图片管理"
$ws.Range("B36").Value = "追加图片"
$ws.Range("B37").Value = "删除图片"

# … then the shared category label (added last so it becomes the final
# new shared-string entry, same order the workbook was originally authored in).
$ws.Range("C35").Value = "商品管理 - 图片管理"
$ws.Range("C36").Value = "商品管理 - 图片管理"
$ws.Range("C37").Value = "商品管理 - 图片管理"

# Move the visible selection down to where the new rows were added.
$ws.Range("D36").Select()
